$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

$ws.Range("B3:B52").Select()
